$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a string value to a cell while forcing Excel to store it
# as literal text (avoids numeric auto-conversion of values like "584.34"
# or "0.0770", which would otherwise lose trailing zeros / formatting).
# NumberFormat is reset back to "Normal" style afterwards so the cell
# keeps the same (default) style it had before the edit.
function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range('D2').Value = '62.928.59'
$ws.Range('E2').Value = '  +2.31%  '
$ws.Range('D3').Value = '3.478.46'
$ws.Range('E3').Value = '  +2.46%  '
$ws.Range('E4').Value = '  +0.29%  '
Set-TextValue $ws.Range('D5') '584.34'
$ws.Range('E5').Value = '  +1.30%  '
Set-TextValue $ws.Range('D6') '147.20'
$ws.Range('E6').Value = '  +4.17%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('E8').Value = '  +0.85%  '
Set-TextValue $ws.Range('D9') '7.69'
$ws.Range('E9').Value = '  -0.22%  '
$ws.Range('E10').Value = '  +2.09%  '
Set-TextValue $ws.Range('D11') '0.397'
$ws.Range('E11').Value = '  +2.30%  '
$ws.Range('D12').Value = '4.086.65'
$ws.Range('E12').Value = '  +2.80%  '
Set-TextValue $ws.Range('D13') '29.75'
$ws.Range('E13').Value = '  +5.35%  '
$ws.Range('E14').Value = '  -0.17%  '
$ws.Range('D15').Value = '3.494.30'
$ws.Range('E15').Value = '  +2.72%  '
$ws.Range('E16').Value = '  +0.78%  '
$ws.Range('D17').Value = '63.124.47'
$ws.Range('E17').Value = '  +2.67%  '
Set-TextValue $ws.Range('D18') '6.30'
$ws.Range('E18').Value = '  +2.58%  '
Set-TextValue $ws.Range('D19') '14.31'
$ws.Range('E19').Value = '  +4.48%  '
Set-TextValue $ws.Range('D20') '9.30'
$ws.Range('E20').Value = '  +3.63%  '
Set-TextValue $ws.Range('D21') '387.34'
$ws.Range('E21').Value = '  -1.09%  '
Set-TextValue $ws.Range('D22') '0.563'
$ws.Range('E22').Value = '  +1.30%  '
Set-TextValue $ws.Range('D23') '75.00'
$ws.Range('E23').Value = '  -0.53%  '
$ws.Range('E24').Value = '  -0.16%  '
$ws.Range('D25').Value = '3.631.92'
$ws.Range('E25').Value = '  +2.86%  '
$ws.Range('E26').Value = '  +2.91%  '
Set-TextValue $ws.Range('D27') '0.180'
$ws.Range('E27').Value = '  -6.45%  '
Set-TextValue $ws.Range('D28') '7.65'
$ws.Range('E28').Value = '  +4.69%  '
$ws.Range('E29').Value = '  +0.09%  '
$ws.Range('E30').Value = '  +2.80%  '
$ws.Range('E31').Value = '  -0.17%  '
$ws.Range('E32').Value = '  +3.96%  '
$ws.Range('E33').Value = '  +0.03%  '
$ws.Range('E34').Value = '  +1.48%  '
Set-TextValue $ws.Range('D35') '5.30'
$ws.Range('E35').Value = '  +4.75%  '
Set-TextValue $ws.Range('D36') '7.09'
$ws.Range('E36').Value = '  +2.12%  '
Set-TextValue $ws.Range('D37') '31.88'
$ws.Range('E37').Value = '  +21.58%  '
$ws.Range('B38').Value = 'Monero'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range('D38') '170.98'
$ws.Range('E38').Value = '  +1.89%  '
$ws.Range('B39').Value = 'ImmutableX'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range('D39') '1.57'
$ws.Range('E39').Value = '  +6.23%  '
$ws.Range('D40').Value = '3.523.25'
$ws.Range('E40').Value = '  +2.80%  '
Set-TextValue $ws.Range('D41') '0.0770'
$ws.Range('E41').Value = '  -0.07%  '
Set-TextValue $ws.Range('D42') '0.807'
$ws.Range('E42').Value = '  +3.58%  '
$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range('D43') '4.49'
$ws.Range('E43').Value = '  +0.94%  '
$ws.Range('B44').Value = 'OKB'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Range('D44') '42.32'
$ws.Range('E44').Value = '  -0.32%  '
$ws.Range('B45').Value = 'ONDO'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
Set-TextValue $ws.Range('D45') '1.21'
$ws.Range('E45').Value = '  +5.84%  '
$ws.Range('B46').Value = 'Stacks'
$ws.Range('C46').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range('D46') '1.71'
$ws.Range('E46').Value = '  +2.66%  '
$ws.Range('D47').Value = '2.619.15'
$ws.Range('E47').Value = '  +5.93%  '
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range('D48') '23.37'
$ws.Range('E48').Value = '  +1.63%  '
$ws.Range('B49').Value = 'dogwifhat'
$ws.Range('C49').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws.Range('D49') '2.26'
$ws.Range('E49').Value = '  +8.69%  '
Set-TextValue $ws.Range('D50') '6.77'
$ws.Range('E50').Value = '  +0.84%  '
$ws.Range('E51').Value = '  +2.37%  '
